# Generate Report for Handback
# Both locales (zh-cn, de-de) are now in sync with en-US: the handback has
# been processed, so the "Ready for handoff" status becomes
# "Handed back: in sync with en-US", the Latest Handback DateTime is
# refreshed, and the stale "handback file is not latest" error is cleared.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Overview tab mirrors each locale's status in columns E (zh-cn) and F (de-de)
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Columns.Item(5).AutoFit()
$ov.Columns.Item(6).AutoFit()

# zh-cn detail tab
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("K2").Value = "2016-09-06 11:05:48"
$zh.Range("P2").Value = ""
$zh.Columns.Item(3).AutoFit()
$zh.Columns.Item(16).AutoFit()

# de-de detail tab
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("K2").Value = "2016-09-06 11:05:56"
$de.Range("P2").Value = ""
$de.Columns.Item(3).AutoFit()
$de.Columns.Item(16).AutoFit()
